$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Append a new log row (row 12) describing the IE meta-tag / Maintenance
# --- Section / Creative Commons update ---
$ws.Range("A12").Value = 42409
$ws.Range("B12").Value = "cbardash"
$ws.Range("C12").Value = "Added meta tag to fix IR problem / Updated Maint Sec layer / added Creative Commones license"
$ws.Range("D12").Value = "Meta tag fixes compatibility view problem with IE.  Added new Maintenance Section layer on AGO and updated service.  Added Creative Commons license statement to top."
$ws.Range("G12").Value = "No"
$ws.Range("E12").Value = "Open in IE with Compatibility View turned on and check that the app looks correct and all functions work.  Turn on Maintenance Section layer and select polygons to test popup."
$ws.Range("F12").Value = "YES"
$ws.Range("I12").Value = "4.1.1 info"

# Match the row height used by the other wrapped/multi-line rows in the log
$ws.Rows.Item(12).RowHeight = 45

# --- Update the view so the newly entered row is in focus, mirroring the
# --- author's selection/scroll position after data entry ---
$win = $excel.ActiveWindow
$win.ScrollRow = 11
$win.ScrollColumn = 4
$null = $ws.Range("E11").Select()
